$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8070710853421019
$ws.Range("C2").Value = 0.1558612813855618
$ws.Range("D2").Value = 0.1434369537791085
$ws.Range("F2").Value = 2.01944559355163
$ws.Range("G2").Value = 0.002515851339151515
$ws.Range("J2").Value = 0.2383690032694119
$ws.Range("K2").Value = 0.3802870271782979
$ws.Range("L2").Value = 0.3084297752315166
$ws.Range("M2").Value = 0.2449838662812631
$ws.Range("N2").Value = 2.524859982403882
$ws.Range("O2").Value = 5.259120794124755
$ws.Range("B3").Value = 0.7750478471484996
$ws.Range("C3").Value = 0.1551571922457455
$ws.Range("D3").Value = 0.1421654971644344
$ws.Range("F3").Value = 2.026840617455591
$ws.Range("G3").Value = 0.002518261474892289
$ws.Range("J3").Value = 0.2397036272473088
$ws.Range("K3").Value = 0.3505862734810279
$ws.Range("L3").Value = 0.3064982543490373
$ws.Range("M3").Value = 0.2392293765846993
$ws.Range("N3").Value = 2.546200127696054
$ws.Range("O3").Value = 5.284425909537447
$ws.Range("B4").Value = 0.7556998502377041
$ws.Range("C4").Value = 0.1547247576522643
$ws.Range("D4").Value = 0.1414334756560578
$ws.Range("F4").Value = 2.032177260091899
$ws.Range("G4").Value = 0.002519821937518898
$ws.Range("J4").Value = 0.2405790339290128
$ws.Range("K4").Value = 0.3324517332338957
$ws.Range("L4").Value = 0.3054202363314857
$ws.Range("M4").Value = 0.2357926733530036
$ws.Range("N4").Value = 2.559972868366446
$ws.Range("O4").Value = 5.30206293338189
$ws.Range("B5").Value = 0.7478951164171974
$ws.Range("C5").Value = 0.1545485215197822
$ws.Range("D5").Value = 0.1411474693906669
$ws.Range("F5").Value = 2.034552446216409
$ws.Range("G5").Value = 0.002520478172698966
$ws.Range("J5").Value = 0.2409498608772687
$ws.Range("K5").Value = 0.3250878530402872
$ws.Range("L5").Value = 0.3050081863469742
$ws.Range("M5").Value = 0.2344166245605273
$ws.Range("N5").Value = 2.565753886708499
$ws.Range("O5").Value = 5.309778697457276
$ws.Range("B6").Value = 0.7466039817411172
$ws.Range("C6").Value = 0.1545192571458216
$ws.Range("D6").Value = 0.1411007229165264
$ws.Range("F6").Value = 2.034958959741516
$ws.Range("G6").Value = 0.002520588369785156
$ws.Range("J6").Value = 0.2410122882297738
$ws.Range("K6").Value = 0.3238666764635099
$ws.Range("L6").Value = 0.304941414831994
$ws.Range("M6").Value = 0.2341896129941716
$ws.Range("N6").Value = 2.566724000858432
$ws.Range("O6").Value = 5.311091833334018
$ws.Range("B7").Value = 0.7555942691434439
$ws.Range("C7").Value = 0.1547223809133627
$ws.Range("D7").Value = 0.1414295685960667
$ws.Range("F7").Value = 2.032208480648407
$ws.Range("G7").Value = 0.002519830705299556
$ws.Range("J7").Value = 0.2405839779417764
$ws.Range("K7").Value = 0.3323523150839947
$ws.Range("L7").Value = 0.3054145687989873
$ws.Range("M7").Value = 0.2357740163386737
$ws.Range("N7").Value = 2.560050150776874
$ws.Range("O7").Value = 5.302164850097753
$ws.Range("B8").Value = 0.795964611501546
$ws.Range("C8").Value = 0.1556185452770009
$ws.Range("D8").Value = 0.1429884937014805
$ws.Range("F8").Value = 2.021830341591595
$ws.Range("G8").Value = 0.002516665654880452
$ws.Range("J8").Value = 0.2388175855877464
$ws.Range("K8").Value = 0.370025350151451
$ws.Range("L8").Value = 0.3077414466239148
$ws.Range("M8").Value = 0.2429797618643867
$ws.Range("N8").Value = 2.532079100547634
$ws.Range("O8").Value = 5.267410572116631
$ws.Range("B9").Value = 0.8776003297137436
$ws.Range("C9").Value = 0.1573743917359991
$ws.Range("D9").Value = 0.1464291137432028
$ws.Range("F9").Value = 2.007783825456748
$ws.Range("G9").Value = 0.002511096066696195
$ws.Range("J9").Value = 0.2357964655439169
$ws.Range("K9").Value = 0.4446927042972391
$ws.Range("L9").Value = 0.3131566715099723
$ws.Range("M9").Value = 0.2578707889298144
$ws.Range("N9").Value = 2.482537856730847
$ws.Range("O9").Value = 5.215893359050568
$ws.Range("B10").Value = 0.9390566953217387
$ws.Range("C10").Value = 0.1586628417522746
$ws.Range("D10").Value = 0.1491877880002619
$ws.Range("F10").Value = 2.001293577559693
$ws.Range("G10").Value = 0.002507388692851829
$ws.Range("J10").Value = 0.2338452426756419
$ws.Range("K10").Value = 0.500014976272297
$ws.Range("L10").Value = 0.317649753015715
$ws.Range("M10").Value = 0.2692685612427539
$ws.Range("N10").Value = 2.449369321575446
$ws.Range("O10").Value = 5.188156704787616
$ws.Range("B11").Value = 0.9673305367142859
$ws.Range("C11").Value = 0.1592485229331331
$ws.Range("D11").Value = 0.1504922979969336
$ws.Range("F11").Value = 1.999169740214327
$ws.Range("G11").Value = 0.002505784824172299
$ws.Range("J11").Value = 0.2330155363563975
$ws.Range("K11").Value = 0.5252798422566798
$ws.Range("L11").Value = 0.31980447340095
$ws.Range("M11").Value = 0.2745517452569715
$ws.Range("N11").Value = 2.434979700460545
$ws.Range("O11").Value = 5.177728731231667
$ws.Range("B12").Value = 0.9780820965864621
$ws.Range("C12").Value = 0.1594702283079386
$ws.Range("D12").Value = 0.150993351960139
$ws.Range("F12").Value = 1.998484395706143
$ws.Range("G12").Value = 0.00250518930301101
$ws.Range("J12").Value = 0.2327096517332699
$ws.Range("K12").Value = 0.5348607320602525
$ws.Range("L12").Value = 0.3206362425561906
$ws.Range("M12").Value = 0.276566346785728
$ws.Range("N12").Value = 2.429631144033989
$ws.Range("O12").Value = 5.174094278155735
$ws.Range("B13").Value = 0.9757645716596812
$ws.Range("C13").Value = 0.1594224838728593
$ws.Range("D13").Value = 0.1508851278629351
$ws.Range("F13").Value = 1.998626712424738
$ws.Range("G13").Value = 0.002505317033830421
$ws.Range("J13").Value = 0.2327751603283481
$ws.Range("K13").Value = 0.5327967180968187
$ws.Range("L13").Value = 0.3204564041736973
$ws.Range("M13").Value = 0.2761318473046828
$ws.Range("N13").Value = 2.430778582766235
$ws.Range("O13").Value = 5.174863047249517
$ws.Range("B14").Value = 0.9682141788796343
$ws.Range("C14").Value = 0.1592667644477572
$ws.Range("D14").Value = 0.1505333787971637
$ws.Range("F14").Value = 1.999110974853579
$ws.Range("G14").Value = 0.00250573559350562
$ws.Range("J14").Value = 0.232990204675751
$ws.Range("K14").Value = 0.5260677976259842
$ws.Range("L14").Value = 0.3198725870504404
$ws.Range("M14").Value = 0.274717208575602
$ws.Range("N14").Value = 2.434537657834696
$ws.Range("O14").Value = 5.177423424310291
$ws.Range("B15").Value = 0.9635951673053285
$ws.Range("C15").Value = 0.1591713709965958
$ws.Range("D15").Value = 0.1503188404094118
$ws.Range("F15").Value = 1.999423077418015
$ws.Range("G15").Value = 0.002505993512385249
$ws.Range("J15").Value = 0.2331230067664265
$ws.Range("K15").Value = 0.5219478993030293
$ws.Range("L15").Value = 0.3195170396364944
$ws.Range("M15").Value = 0.273852516691484
$ws.Range("N15").Value = 2.436853284819776
$ws.Range("O15").Value = 5.179032658869659
$ws.Range("B16").Value = 0.9372152513352319
$ws.Range("C16").Value = 0.1586245558689896
$ws.Range("D16").Value = 0.1491035272683803
$ws.Range("F16").Value = 2.001449027988031
$ws.Range("G16").Value = 0.002507495167754541
$ws.Range("J16").Value = 0.2339006295365333
$ws.Range("K16").Value = 0.4983657944709989
$ws.Range("L16").Value = 0.3175111571420075
$ws.Range("M16").Value = 0.2689252568720875
$ws.Range("N16").Value = 2.450323779114392
$ws.Range("O16").Value = 5.188882217488754
$ws.Range("B17").Value = 0.9211127154293308
$ws.Range("C17").Value = 0.1582889781412788
$ws.Range("D17").Value = 0.1483706235090239
$ws.Range("F17").Value = 2.002903934272481
$ws.Range("G17").Value = 0.002508437510977045
$ws.Range("J17").Value = 0.2343924943423037
$ws.Range("K17").Value = 0.4839237856963052
$ws.Range("L17").Value = 0.3163089129844394
$ws.Range("M17").Value = 0.2659276019643571
$ws.Range("N17").Value = 2.458766482676479
$ws.Range("O17").Value = 5.195485073900528
$ws.Range("B18").Value = 0.9118808462821164
$ws.Range("C18").Value = 0.1580959221309968
$ws.Range("D18").Value = 0.1479537467630649
$ws.Range("F18").Value = 2.003818765186722
$ws.Range("G18").Value = 0.00250898730296276
$ws.Range("J18").Value = 0.2346808544209171
$ws.Range("K18").Value = 0.475626425941158
$ws.Range("L18").Value = 0.3156278490064608
$ws.Range("M18").Value = 0.2642126902055892
$ws.Range("N18").Value = 2.463688288478298
$ws.Range("O18").Value = 5.199488982509934
$ws.Range("B19").Value = 0.908760248457213
$ws.Range("C19").Value = 0.1580305501995554
$ws.Range("D19").Value = 0.147813403522008
$ws.Range("F19").Value = 2.004141916558808
$ws.Range("G19").Value = 0.002509174791326731
$ws.Range("J19").Value = 0.2347794252928086
$ws.Range("K19").Value = 0.4728186961823155
$ws.Range("L19").Value = 0.3153990479967703
$ws.Range("M19").Value = 0.263633645793
$ws.Range("N19").Value = 2.465366026882759
$ws.Range("O19").Value = 5.200880052663962
$ws.Range("B20").Value = 0.9228237702313038
$ws.Range("C20").Value = 0.1583247052778631
$ws.Range("D20").Value = 0.1484481594942082
$ws.Range("F20").Value = 2.002740985337013
$ws.Range("G20").Value = 0.002508336392134736
$ws.Range("J20").Value = 0.2343395703381077
$ws.Range("K20").Value = 0.4854602033219919
$ws.Range("L20").Value = 0.3164358147268445
$ws.Range("M20").Value = 0.2662457502110698
$ws.Range("N20").Value = 2.457860934029377
$ws.Range("O20").Value = 5.194760859249016
$ws.Range("B21").Value = 0.9704307002437531
$ws.Range("C21").Value = 0.1593125052985158
$ws.Range("D21").Value = 0.1506365048147558
$ws.Range("F21").Value = 1.998965510110267
$ws.Range("G21").Value = 0.002505612331651365
$ws.Range("J21").Value = 0.2329268156507958
$ws.Range("K21").Value = 0.5280438788920208
$ws.Range("L21").Value = 0.3200436397191453
$ws.Range("M21").Value = 0.2751323442613227
$ws.Range("N21").Value = 2.43343079950599
$ws.Range("O21").Value = 5.176662850619863
$ws.Range("B22").Value = 1.0018057822966
$ws.Range("C22").Value = 0.1599576196979129
$ws.Range("D22").Value = 0.152107861445387
$ws.Range("F22").Value = 1.997190983391789
$ws.Range("G22").Value = 0.002503900925642864
$ws.Range("J22").Value = 0.2320519096421307
$ws.Range("K22").Value = 0.5559539924161925
$ws.Range("L22").Value = 0.3224937417764693
$ws.Range("M22").Value = 0.2810216272547734
$ws.Range("N22").Value = 2.418049965169027
$ws.Range("O22").Value = 5.166667071770689
$ws.Range("B23").Value = 0.9850366483290429
$ws.Range("C23").Value = 0.1596133582582979
$ws.Range("D23").Value = 0.1513188271240153
$ws.Range("F23").Value = 1.998074756698969
$ws.Range("G23").Value = 0.002504808046248622
$ws.Range("J23").Value = 0.2325144407196227
$ws.Range("K23").Value = 0.5410507658556298
$ws.Range("L23").Value = 0.3211776776961699
$ws.Range("M23").Value = 0.2778710120815546
$ws.Range("N23").Value = 2.426205423495233
$ws.Range("O23").Value = 5.171834502718923
$ws.Range("B24").Value = 0.92205012226907
$ws.Range("C24").Value = 0.158308553438161
$ws.Range("D24").Value = 0.1484130915182504
$ws.Range("F24").Value = 2.002814410365588
$ws.Range("G24").Value = 0.002508382083013897
$ws.Range("J24").Value = 0.2343634798962899
$ws.Range("K24").Value = 0.4847655716066299
$ws.Range("L24").Value = 0.3163784109119661
$ws.Range("M24").Value = 0.2661018889628366
$ws.Range("N24").Value = 2.458270120805663
$ws.Range("O24").Value = 5.19508762925264
$ws.Range("B25").Value = 0.8552541002107432
$ws.Range("C25").Value = 0.1568996186500158
$ws.Range("D25").Value = 0.1454575795724082
$ws.Range("F25").Value = 2.010910233654243
$ws.Range("G25").Value = 0.002512534977040972
$ws.Range("J25").Value = 0.2365665157960155
$ws.Range("K25").Value = 0.4244104042804508
$ws.Range("L25").Value = 0.3116009927686747
$ws.Range("M25").Value = 0.2537616009857437
$ws.Range("N25").Value = 2.495372415747093
$ws.Range("O25").Value = 5.228052115431808
